# Add Port-Channel Support to L2
# Reworks the lab-switch sheet so the "description/name" column moves from D
# to G, freeing up D/E/F for "allowed vlans" / "channel-group" / "mode", and
# adds two port-channel rows (17, 18) describing the new LAG to the Cisco GW.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Style donor: every populated cell in this sheet uses cell style index 1
# (s="1"). Keep a single cell on the clipboard and stamp-copy its format
# (not its value) onto every brand-new cell before setting the real value.
$styleDonor = $ws.Cells.Item(1, 1)

function Stamp-Style($cell) {
    $styleDonor.Copy() | Out-Null
    $cell.PasteSpecial(-4122) | Out-Null
}

# ---------------------------------------------------------------------
# 1. Header row: insert "allowed vlans" (D, replacing the old
#    "description/name" text), then the brand new "channel-group" (E),
#    "mode" (F) and, finally, "description/name" again in its new home (G).
# ---------------------------------------------------------------------
$ws.Cells.Item(1, 4).Value = "allowed vlans"

foreach ($c in 5..7) {
    Stamp-Style $ws.Cells.Item(1, $c)
}
$ws.Cells.Item(1, 5).Value = "channel-group"
$ws.Cells.Item(1, 6).Value = "mode"
$ws.Cells.Item(1, 7).Value = "description/name"

# ---------------------------------------------------------------------
# 2. Rows 2-8 and 13-16: the "description/name" text currently sitting in
#    column D moves, unchanged, to the new column G; column D is then
#    fully removed (not merely blanked) for these rows.
# ---------------------------------------------------------------------
$moveRows = @(2, 3, 4, 5, 6, 7, 8, 13, 14, 15, 16)
foreach ($r in $moveRows) {
    $dCell = $ws.Cells.Item($r, 4)
    $text = $dCell.Value2
    $gCell = $ws.Cells.Item($r, 7)
    Stamp-Style $gCell
    $gCell.Value = $text
    $dCell.Clear() | Out-Null
}

# ---------------------------------------------------------------------
# 3. Rows 9-12 (fastethernet 0/8 .. 0/11): these had only interface
#    type/number before. They now also carry "trunk" (C), a channel-group
#    number (E) and a mode (F, "on" for the first pair, "active" for the
#    second).
# ---------------------------------------------------------------------
$channelRows = @(
    @{ row = 9;  group = 100; mode = "on" },
    @{ row = 10; group = 100; mode = "on" },
    @{ row = 11; group = 101; mode = "active" },
    @{ row = 12; group = 101; mode = "active" }
)
foreach ($item in $channelRows) {
    $r = $item.row

    $cCell = $ws.Cells.Item($r, 3)
    Stamp-Style $cCell
    $cCell.Value = "trunk"

    $eCell = $ws.Cells.Item($r, 5)
    Stamp-Style $eCell
    $eCell.Value = $item.group

    $fCell = $ws.Cells.Item($r, 6)
    Stamp-Style $fCell
    $fCell.Value = $item.mode
}

# ---------------------------------------------------------------------
# 4. Row 17 (already present, but empty aside from styled A17/C17):
#    becomes the first Port-Channel row, and row 18 is a brand new second
#    Port-Channel row. Both describe the LAG to "Agg cisco-GW", leaving F
#    (mode) intentionally blank - but still styled - since port-channel
#    interfaces don't have a channel-group mode of their own.
# ---------------------------------------------------------------------
$portChannelRows = @(
    @{ row = 17; num = 100 },
    @{ row = 18; num = 101 }
)
foreach ($item in $portChannelRows) {
    $r = $item.row

    $aCell = $ws.Cells.Item($r, 1)
    Stamp-Style $aCell
    $aCell.Value = "port-channel"

    $bCell = $ws.Cells.Item($r, 2)
    Stamp-Style $bCell
    $bCell.Value = $item.num

    $cCell = $ws.Cells.Item($r, 3)
    Stamp-Style $cCell
    $cCell.Value = "trunk"

    $fCell = $ws.Cells.Item($r, 6)
    Stamp-Style $fCell

    $gCell = $ws.Cells.Item($r, 7)
    Stamp-Style $gCell
    $gCell.Value = "Agg cisco-GW"
}

# ---------------------------------------------------------------------
# 5. Column widths for the new/changed layout (B..G).
#
# Excel's ColumnWidth setter snaps to whole pixels (width_in_chars ->
# round(chars * MaxDigitWidth) + 5px -> stored width = px / MaxDigitWidth,
# MaxDigitWidth = 6 for this workbook's Normal font), so it cannot hit an
# arbitrary two-decimal "width" value exactly. The inputs below are chosen
# so the resulting stored width lands on the closest achievable pixel to
# the target (8.0, 5.14, 11.86, 12.57, 5.71, 15.14).
# ---------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 7.166666666666667
$ws.Columns.Item(3).ColumnWidth = 4.333333333333333
$ws.Columns.Item(4).ColumnWidth = 11.0
$ws.Columns.Item(5).ColumnWidth = 11.666666666666666
$ws.Columns.Item(6).ColumnWidth = 4.833333333333333
$ws.Columns.Item(7).ColumnWidth = 14.333333333333334

Write-Host "Port-channel support added"
